# Edit script: add new wave "12. 10. 2021" / "20. 10. 2021" update columns
# to both sheets of the ZBP_11_obavy_epidemie workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "data": new column AJ (percentages), header date "12. 10. 2021"
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data")

# Copy formatting of the previous header cell (AI1, style "1") onto the
# new header cell AJ1, then set its text.
$wsData.Range("AI1").Copy() | Out-Null
$wsData.Range("AJ1").PasteSpecial(-4122) | Out-Null
$wsData.Range("AJ1").Value = "12. 10. 2021"

# New data values for AJ2:AJ76 (rows 2-76, one per row, in order).
$dataValues = 0.28,0.45,0.27,0.34,0.49,0.17,0.31,0.47,0.22,0.22,0.4,0.38,0.27,0.41,0.32,0.26,0.43,0.31,0.29,0.46,0.25,0.35,0.42,0.23,0.22,0.47,0.31,0.28,0.45,0.27,0.29,0.45,0.26,0.34,0.38,0.28,0.26,0.47,0.27,0.2,0.48,0.32,0.31,0.44,0.25,0.33,0.43,0.24,0.49,0.21,0.3,0.32,0.46,0.22,0.29,0.48,0.23,0.42,0.44,0.14,0.38,0.42,0.2,0.26,0.42,0.32,0.26,0.55,0.19,0.24,0.52,0.24,0.2,0.38,0.42
for ($i = 0; $i -lt $dataValues.Length; $i++) {
    $wsData.Cells.Item($i + 2, 36).Value = $dataValues[$i]
}

# Update the footnote text in row 77 with the new "aktualizace" date.
$wsData.Range("A77").Value = "Život během pandemie, Obavy z epidemie, % respondentů celkově a ve skupinách, aktualizace 20. 10. 2021"

# ---------------------------------------------------------------------
# Sheet "pocetR": new column AI (sample sizes), header date "12. 10. 2021"
# ---------------------------------------------------------------------
$wsPocet = $wb.Worksheets.Item("pocetR")

# Copy formatting of the previous header cell (AH1, style "2") onto the
# new header cell AI1, then set its text.
$wsPocet.Range("AH1").Copy() | Out-Null
$wsPocet.Range("AI1").PasteSpecial(-4122) | Out-Null
$wsPocet.Range("AI1").Value = "12. 10. 2021"

# New sample-size values for AI2:AI26 (rows 2-26, one per row, in order).
$pocetValues = 1836,454,670,712,296,316,1224,895,941,960,418,216,242,44,144,92,20,258,492,239,336,319,215,334,393
for ($i = 0; $i -lt $pocetValues.Length; $i++) {
    $wsPocet.Cells.Item($i + 2, 35).Value = $pocetValues[$i]
}

# Update the footnote text in row 27 with the new "aktualizace" date.
$wsPocet.Range("A27").Value = "Život během pandemie, Obavy z epidemie, velikost dotázaného souboru celkově a ve skupinách, aktualizace 20. 10. 2021"

# Extend the formatting of the blank footnote row to the new column AI27,
# matching the existing blank-but-present cells B27:AH27.
$wsPocet.Range("AI27").Style = "Normal"
